$wb = $excel.ActiveWorkbook

# --- "Date Formula test" sheet: add a C column with a "timestamp" (NOW()) sample ---
$ws = $wb.Worksheets.Item("Date Formula test")

# Give the new C6 cell a date+time number format *before* putting the volatile
# formula in, so the engine doesn't invent its own custom format for it.
$ws.Range("C6").NumberFormat = "m/d/yy h:mm"
$ws.Range("C6").Formula = "=NOW()"

# New header cell C5 ("timestamp") should look like the existing header B5
# (same grey fill style), so copy B5's formatting over before setting the text.
$ws.Range("B5").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "timestamp"

# Widen the new column roughly to fit its content.
$ws.Columns.Item(3).ColumnWidth = 15.83

# Make this the active sheet/tab with cell D10 selected, like in the edited file.
$ws.Activate()
$ws.Range("D10").Select()
